$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 21; this shifts existing rows 21-84 down to 22-85
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with a new weekly record.
# Static/template columns copied from the (shifted) row below (now row 22),
# which previously occupied row 21.
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44607
$ws.Range("D21").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 100112030
$ws.Range("G21").Value = "Poroto granado"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 60
$ws.Range("K21").Value = 22000
$ws.Range("L21").Value = 23000
$ws.Range("M21").Value = 22500
$ws.Range("N21").Value = "$/saco 25 kilos"
$ws.Range("O21").Value = "Provincia de Diguillín"
$ws.Range("P21").Value = 900
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
